$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
# Ensure target cells are treated as text so numeric-looking strings
# (e.g. "13.80", "1.00") keep their exact formatting instead of being
# auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.354.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.457.08'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.94'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.455.89'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.133'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.53%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.416'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.054.12'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.61'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.224.74'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000173'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.450.86'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.97'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.80'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '373.87'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.18%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.28'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000125'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.82'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.32%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.86'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.61%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.71'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.88%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.28'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -7.11%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.73'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.08%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '28.47'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.82'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.65'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.776.25'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.04%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0693'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.28'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '338.65'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '39.97'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0293'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.10%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.995'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.26%  '
